# Weekly update: insert a new price record as row 9 ("Fruta / hortaliza, semanal").
# Every existing data row from 9 downward shifts down by one (old row 79 -> row 80),
# and the newly inserted row 9 carries the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 9..79 down to 10..80, creating a blank row 9.
$ws.Rows.Item(9).EntireRow.Insert()

# Populate the new row 9 with this week's record.
$ws.Cells.Item(9, 1).Value = 2
$ws.Cells.Item(9, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(9, 3).Value = "Coquimbo"
$ws.Cells.Item(9, 4).Value = 44882
$ws.Cells.Item(9, 5).Value = 4
$ws.Cells.Item(9, 6).Value = 100112026
$ws.Cells.Item(9, 7).Value = "Haba"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 560
$ws.Cells.Item(9, 11).Value = 4500
$ws.Cells.Item(9, 12).Value = 5000
$ws.Cells.Item(9, 13).Value = 4750
$ws.Cells.Item(9, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(9, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 16).Value = 190
$ws.Cells.Item(9, 17).Value = 25
$ws.Cells.Item(9, 18).Value = "Hortaliza"
